$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 3
$ws.Range("B2").Value = 160

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 142

$ws.Range("A4").Value = 0
$ws.Range("B4").Value = 139

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 137

$ws.Range("A6").Value = 1
$ws.Range("B6").Value = 124

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 63

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 56
